# Add a new "NutrShortDesc" column (F) to the nutr_def sheet giving a short,
# human-friendly description for each nutrient.
#
# The rows are written in a specific order (not top-to-bottom) so that the
# newly introduced shared-string entries land in the same sequence as in the
# target workbook (this mirrors how the original data was produced
# programmatically rather than typed row-by-row in the UI).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$order = @(
    @(1,  "NutrShortDesc"),
    @(4,  "Carbohydrate"),
    @(11, "Folate"),
    @(12, "Iron"),
    @(13, "Magnesium"),
    @(15, "Phosphorus"),
    @(16, "Potassium"),
    @(22, "Vitamin B12"),
    @(23, "Vitamin B6"),
    @(24, "Vitamin C"),
    @(25, "Vitamin D"),
    @(26, "Vitamin E"),
    @(27, "Vitamin K"),
    @(28, "Zinc"),
    @(29, "Sugar"),
    @(31, "Alcohol"),
    @(3,  "Calcium"),
    @(21, "Vitamin A"),
    @(2,  "Caffeine"),
    @(5,  "Cholesterol"),
    @(6,  "Energy"),
    @(7,  "Fatty acids, total monounsaturated"),
    @(8,  "Fatty acids, total polyunsaturated"),
    @(9,  "Fatty acids, total saturated"),
    @(10, "Fiber, total dietary"),
    @(14, "Niacin"),
    @(17, "Protein"),
    @(18, "Riboflavin"),
    @(19, "Thiamin"),
    @(20, "Total lipid (fat)"),
    @(30, "Fatty acids, total trans")
)

foreach ($pair in $order) {
    $r = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($r, 6).Value = $val
}

# The last row's new cell picks up the same distinctive font style already
# used by D31 ("Alcohol, ethyl"). Copy/PasteSpecial formats reuses the
# existing style instead of fabricating a new (duplicate) one.
$ws.Cells.Item(31, 4).Copy()
$ws.Cells.Item(31, 6).PasteSpecial(-4122)

# Match the saved selection state recorded in the target workbook.
$ws.Range("L18").Select()
